$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formatting from F/G (old D/E) onto new D/E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set values for every data cell D:M per row to match the refreshed financials
# Row 7
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
# Row 8
$ws.Range("D8").Value = 8334100
$ws.Range("E8").Value = 6754400
$ws.Range("F8").Value = 6730800
$ws.Range("G8").Value = 5902800
$ws.Range("H8").Value = 7635100
$ws.Range("I8").Value = 6271900
$ws.Range("J8").Value = 6486200
$ws.Range("K8").Value = 6203600
$ws.Range("L8").Value = 7891500
$ws.Range("M8").Value = 6309200
# Row 9
$ws.Range("D9").Value = 2177800
$ws.Range("E9").Value = 2061100
$ws.Range("F9").Value = 2036400
$ws.Range("G9").Value = 1790700
$ws.Range("H9").Value = 1963500
$ws.Range("I9").Value = 1794100
$ws.Range("J9").Value = 1920800
$ws.Range("K9").Value = 1966100
$ws.Range("L9").Value = 2021300
$ws.Range("M9").Value = 1817100
# Row 10
$ws.Range("D10").Value = 6156400
$ws.Range("E10").Value = 4693300
$ws.Range("F10").Value = 4694400
$ws.Range("G10").Value = 4112100
$ws.Range("H10").Value = 5671700
$ws.Range("I10").Value = 4477900
$ws.Range("J10").Value = 4565400
$ws.Range("K10").Value = 4237500
$ws.Range("L10").Value = 5870200
$ws.Range("M10").Value = 4492200
# Row 12
$ws.Range("D12").Value = 1054700
$ws.Range("E12").Value = 1027700
$ws.Range("F12").Value = 1061400
$ws.Range("G12").Value = 911100
$ws.Range("H12").Value = 980600
$ws.Range("I12").Value = 872900
$ws.Range("J12").Value = 942500
$ws.Range("K12").Value = 996600
$ws.Range("L12").Value = 1006000
$ws.Range("M12").Value = 895600
# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
# Row 14
$ws.Range("D14").Value = 164900
$ws.Range("E14").Value = 164900
$ws.Range("F14").Value = 179500
$ws.Range("G14").Value = 154800
$ws.Range("H14").Value = 164900
$ws.Range("I14").Value = 80800
$ws.Range("J14").Value = 435300
$ws.Range("K14").Value = 189000
$ws.Range("L14").Value = 213600
$ws.Range("M14").Value = 193700
# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
# Row 17
$ws.Range("D17").Value = 5642500
$ws.Range("E17").Value = 5367600
$ws.Range("F17").Value = 5559500
$ws.Range("G17").Value = 4752700
$ws.Range("H17").Value = 5431600
$ws.Range("I17").Value = 4797600
$ws.Range("J17").Value = 5447300
$ws.Range("K17").Value = 5413600
$ws.Range("L17").Value = 5601400
$ws.Range("M17").Value = 5015700
# Row 18
$ws.Range("D18").Value = 2691700
$ws.Range("E18").Value = 1386800
$ws.Range("F18").Value = 1171400
$ws.Range("G18").Value = 1150000
$ws.Range("H18").Value = 2203600
$ws.Range("I18").Value = 1474300
$ws.Range("J18").Value = 1039000
$ws.Range("K18").Value = 790000
$ws.Range("L18").Value = 2290100
$ws.Range("M18").Value = 1293500
# Row 20
$ws.Range("D20").Value = -107700
$ws.Range("E20").Value = 47100
$ws.Range("F20").Value = -16800
$ws.Range("G20").Value = -39300
$ws.Range("H20").Value = 111100
$ws.Range("I20").Value = 85300
$ws.Range("J20").Value = -20200
$ws.Range("K20").Value = -5900
$ws.Range("L20").Value = 14100
$ws.Range("M20").Value = -104500
# Row 21
$ws.Range("D21").Value = 3005800
$ws.Range("E21").Value = 1827700
$ws.Range("F21").Value = 1527000
$ws.Range("G21").Value = 1450700
$ws.Range("H21").Value = 2663600
$ws.Range("I21").Value = 1917500
$ws.Range("J21").Value = 1380000
$ws.Range("K21").Value = 1159700
$ws.Range("L21").Value = 2698600
$ws.Range("M21").Value = 1561200
# Row 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
# Row 23
$ws.Range("D23").Value = 2583900
$ws.Range("E23").Value = 1433900
$ws.Range("F23").Value = 1154500
$ws.Range("G23").Value = 1110800
$ws.Range("H23").Value = 2314700
$ws.Range("I23").Value = 1559600
$ws.Range("J23").Value = 1018800
$ws.Range("K23").Value = 784100
$ws.Range("L23").Value = 2304200
$ws.Range("M23").Value = 1189100
# Row 24
$ws.Range("D24").Value = 691100
$ws.Range("E24").Value = 341100
$ws.Range("F24").Value = 346700
$ws.Range("G24").Value = 316400
$ws.Range("H24").Value = 231100
$ws.Range("I24").Value = 446600
$ws.Range("J24").Value = 270400
$ws.Range("K24").Value = 162000
$ws.Range("L24").Value = 513000
$ws.Range("M24").Value = 338100
# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
# Row 26
$ws.Range("D26").Value = 1892800
$ws.Range("E26").Value = 1092800
$ws.Range("F26").Value = 807800
$ws.Range("G26").Value = 794400
$ws.Range("H26").Value = 2083500
$ws.Range("I26").Value = 1113000
$ws.Range("J26").Value = 748400
$ws.Range("K26").Value = 622100
$ws.Range("L26").Value = 1791200
$ws.Range("M26").Value = 851000
# Row 27
$ws.Range("D27").Value = 1888300
$ws.Range("E27").Value = 1091700
$ws.Range("F27").Value = 806700
$ws.Range("G27").Value = 794400
$ws.Range("H27").Value = 2060000
$ws.Range("I27").Value = 1102900
$ws.Range("J27").Value = 749500
$ws.Range("K27").Value = 611600
$ws.Range("L27").Value = 1793600
$ws.Range("M27").Value = 856900
# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
# Row 29
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"
$ws.Range("M29").Value = "NA"
# Row 30
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
# Row 31
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
# Row 32
$ws.Range("D32").Value = 107700
$ws.Range("E32").Value = -47100
$ws.Range("F32").Value = 16800
$ws.Range("G32").Value = 39300
$ws.Range("H32").Value = -111100
$ws.Range("I32").Value = -85300
$ws.Range("J32").Value = 20200
$ws.Range("K32").Value = 5900
$ws.Range("L32").Value = -14100
$ws.Range("M32").Value = 104500
# Row 33
$ws.Range("D33").Value = 1888300
$ws.Range("E33").Value = 1091700
$ws.Range("F33").Value = 806700
$ws.Range("G33").Value = 794400
$ws.Range("H33").Value = 2060000
$ws.Range("I33").Value = 1102900
$ws.Range("J33").Value = 749500
$ws.Range("K33").Value = 611600
$ws.Range("L33").Value = 1793600
$ws.Range("M33").Value = 856900
# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
# Row 35
$ws.Range("D35").Value = 1888300
$ws.Range("E35").Value = 1091700
$ws.Range("F35").Value = 806700
$ws.Range("G35").Value = 794400
$ws.Range("H35").Value = 2060000
$ws.Range("I35").Value = 1102900
$ws.Range("J35").Value = 749500
$ws.Range("K35").Value = 611600
$ws.Range("L35").Value = 1793600
$ws.Range("M35").Value = 856900
# Row 38
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
# Row 41
$ws.Range("D41").Value = 9679400
$ws.Range("E41").Value = 5056800
$ws.Range("F41").Value = 5065800
$ws.Range("G41").Value = 8524900
$ws.Range("H41").Value = 4500300
$ws.Range("I41").Value = 4734800
$ws.Range("J41").Value = 4752700
$ws.Range("K41").Value = 6968900
$ws.Range("L41").Value = 4345400
$ws.Range("M41").Value = 4826700
# Row 42
$ws.Range("D42").Value = 502700
$ws.Range("E42").Value = 543000
$ws.Range("F42").Value = 444300
$ws.Range("G42").Value = 962700
$ws.Range("H42").Value = 1154500
$ws.Range("I42").Value = 1136600
$ws.Range("J42").Value = 973900
$ws.Range("K42").Value = 1813500
$ws.Range("L42").Value = 1319400
$ws.Range("M42").Value = 588100
# Row 43
$ws.Range("D43").Value = 7466800
$ws.Range("E43").Value = 6030700
$ws.Range("F43").Value = 6205700
$ws.Range("G43").Value = 6478400
$ws.Range("H43").Value = 20987900
$ws.Range("I43").Value = 5868000
$ws.Range("J43").Value = 6488500
$ws.Range("K43").Value = 345100
$ws.Range("L43").Value = 273500
$ws.Range("M43").Value = 352100
# Row 44
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
# Row 45
$ws.Range("D45").Value = 998600
$ws.Range("E45").Value = 1150000
$ws.Range("F45").Value = 1128700
$ws.Range("G45").Value = 952600
$ws.Range("H45").Value = 989600
$ws.Range("I45").Value = 848200
$ws.Range("J45").Value = 842600
$ws.Range("K45").Value = 909700
$ws.Range("L45").Value = 682000
$ws.Range("M45").Value = 747700
# Row 46
$ws.Range("D46").Value = 18647500
$ws.Range("E46").Value = 12780600
$ws.Range("F46").Value = 12844500
$ws.Range("G46").Value = 16918500
$ws.Range("H46").Value = 13385300
$ws.Range("I46").Value = 12587600
$ws.Range("J46").Value = 13057700
$ws.Range("K46").Value = 17779700
$ws.Range("L46").Value = 13573900
$ws.Range("M46").Value = 12177100
# Row 47
$ws.Range("D47").Value = 2301200
$ws.Range("E47").Value = 2268700
$ws.Range("F47").Value = 2207000
$ws.Range("G47").Value = 1919700
$ws.Range("H47").Value = 1955600
$ws.Range("I47").Value = 2009500
$ws.Range("J47").Value = 2305700
$ws.Range("K47").Value = 2371100
$ws.Range("L47").Value = 2270100
$ws.Range("M47").Value = 2245500
# Row 48
$ws.Range("D48").Value = 3986400
$ws.Range("E48").Value = 3837200
$ws.Range("F48").Value = 3684600
$ws.Range("G48").Value = 3415300
$ws.Range("H48").Value = 3328900
$ws.Range("I48").Value = 3148300
$ws.Range("J48").Value = 3050700
$ws.Range("K48").Value = 3075400
$ws.Range("L48").Value = 3028400
$ws.Range("M48").Value = 2785500
# Row 49
$ws.Range("D49").Value = 30239900
$ws.Range("E49").Value = 30033400
$ws.Range("F49").Value = 30024500
$ws.Range("G49").Value = 26582200
$ws.Range("H49").Value = 27194800
$ws.Range("I49").Value = 27366500
$ws.Range("J49").Value = 28298800
$ws.Range("K49").Value = 31375900
$ws.Range("L49").Value = 31806700
$ws.Range("M49").Value = 30526100
# Row 50
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
# Row 52
$ws.Range("D52").Value = 2597400
$ws.Range("E52").Value = 2294500
$ws.Range("F52").Value = 2268700
$ws.Range("G52").Value = 2184500
$ws.Range("H52").Value = 1935400
$ws.Range("I52").Value = 1372200
$ws.Range("J52").Value = 1420400
$ws.Range("K52").Value = 1416800
$ws.Range("L52").Value = 1293500
$ws.Range("M52").Value = 1097500
# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
# Row 54
$ws.Range("D54").Value = 57772400
$ws.Range("E54").Value = 51214400
$ws.Range("F54").Value = 51029200
$ws.Range("G54").Value = 51020300
$ws.Range("H54").Value = 47666600
$ws.Range("I54").Value = 46484000
$ws.Range("J54").Value = 48133400
$ws.Range("K54").Value = 56018900
$ws.Range("L54").Value = 51972800
$ws.Range("M54").Value = 48831700
# Row 57
$ws.Range("D57").Value = 1667300
$ws.Range("E57").Value = 1398000
$ws.Range("F57").Value = 1318300
$ws.Range("G57").Value = 1196000
$ws.Range("H57").Value = 1291400
$ws.Range("I57").Value = 1221800
$ws.Range("J57").Value = 1281300
$ws.Range("K57").Value = 1388600
$ws.Range("L57").Value = 1503700
$ws.Range("M57").Value = 1307600
# Row 58
$ws.Range("D58").Value = 861700
$ws.Range("E58").Value = 1260000
$ws.Range("F58").Value = 1309400
$ws.Range("G58").Value = 1788500
$ws.Range("H58").Value = 3207800
$ws.Range("I58").Value = 1446200
$ws.Range("J58").Value = 608100
$ws.Range("K58").Value = 2083500
$ws.Range("L58").Value = 2128100
$ws.Range("M58").Value = 1489600
# Row 59
$ws.Range("D59").Value = 9230600
$ws.Range("E59").Value = 8667400
$ws.Range("F59").Value = 10204500
$ws.Range("G59").Value = 9839900
$ws.Range("H59").Value = 12110800
$ws.Range("I59").Value = 7998700
$ws.Range("J59").Value = 9687300
$ws.Range("K59").Value = 11375400
$ws.Range("L59").Value = 7723700
$ws.Range("M59").Value = 7707200
# Row 60
$ws.Range("D60").Value = 11759600
$ws.Range("E60").Value = 11325400
$ws.Range("F60").Value = 12832200
$ws.Range("G60").Value = 12824300
$ws.Range("H60").Value = 11455500
$ws.Range("I60").Value = 10666800
$ws.Range("J60").Value = 11576700
$ws.Range("K60").Value = 14847500
$ws.Range("L60").Value = 11355400
$ws.Range("M60").Value = 10504400
# Row 61
$ws.Range("D61").Value = 11820200
$ws.Range("E61").Value = 7631800
$ws.Range("F61").Value = 7303000
$ws.Range("G61").Value = 7269400
$ws.Range("H61").Value = 5612200
$ws.Range("I61").Value = 6519900
$ws.Range("J61").Value = 6984400
$ws.Range("K61").Value = 7575800
$ws.Range("L61").Value = 7569900
$ws.Range("M61").Value = 8507800
# Row 62
$ws.Range("D62").Value = 1792900
$ws.Range("E62").Value = 1956800
$ws.Range("F62").Value = 1882700
$ws.Range("G62").Value = 2104900
$ws.Range("H62").Value = 2093600
$ws.Range("I62").Value = 1954500
$ws.Range("J62").Value = 2055500
$ws.Range("K62").Value = 2175100
$ws.Range("L62").Value = 2062400
$ws.Range("M62").Value = 1925000
# Row 63
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
# Row 64
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
# Row 65
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
# Row 66
$ws.Range("D66").Value = 25423200
$ws.Range("E66").Value = 20962100
$ws.Range("F66").Value = 22067300
$ws.Range("G66").Value = 22232200
$ws.Range("H66").Value = 19073800
$ws.Range("I66").Value = 19142300
$ws.Range("J66").Value = 20624400
$ws.Range("K66").Value = 24634800
$ws.Range("L66").Value = 21012400
$ws.Range("M66").Value = 20963100
# Row 68
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
# Row 69
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
# Row 70
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
# Row 71
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
# Row 72
$ws.Range("D72").Value = 30750400
$ws.Range("E72").Value = 30039000
$ws.Range("F72").Value = 27756900
$ws.Range("G72").Value = 28551300
$ws.Range("H72").Value = 28360500
$ws.Range("I72").Value = 26873900
$ws.Range("J72").Value = 24688300
$ws.Range("K72").Value = 30535500
$ws.Range("L72").Value = 30105900
$ws.Range("M72").Value = 27054000
# Row 73
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
# Row 74
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
# Row 75
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
# Row 76
$ws.Range("D76").Value = 32349200
$ws.Range("E76").Value = 30252200
$ws.Range("F76").Value = 28961900
$ws.Range("G76").Value = 28788000
$ws.Range("H76").Value = 28592800
$ws.Range("I76").Value = 27341800
$ws.Range("J76").Value = 27509000
$ws.Range("K76").Value = 31384200
$ws.Range("L76").Value = 30960400
$ws.Range("M76").Value = 27868600
# Row 77
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
# Row 80
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
# Row 81
$ws.Range("D81").Value = 1888300
$ws.Range("E81").Value = 1091700
$ws.Range("F81").Value = 806700
$ws.Range("G81").Value = 794400
$ws.Range("H81").Value = 2060000
$ws.Range("I81").Value = 1102900
$ws.Range("J81").Value = 749500
$ws.Range("K81").Value = 611600
$ws.Range("L81").Value = 1793600
$ws.Range("M81").Value = 856900
# Row 83
$ws.Range("D83").Value = 421900
$ws.Range("E83").Value = 393800
$ws.Range("F83").Value = 372500
$ws.Range("G83").Value = 340000
$ws.Range("H83").Value = 348900
$ws.Range("I83").Value = 357900
$ws.Range("J83").Value = 361300
$ws.Range("K83").Value = 375600
$ws.Range("L83").Value = 394400
$ws.Range("M83").Value = 372100
# Row 84
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
# Row 85
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
# Row 86
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
# Row 87
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
# Row 88
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
# Row 89
$ws.Range("D89").Value = 918900
$ws.Range("E89").Value = 559900
$ws.Range("F89").Value = 457800
$ws.Range("G89").Value = 2891400
$ws.Range("H89").Value = 1032200
$ws.Range("I89").Value = 686700
$ws.Range("J89").Value = 719200
$ws.Range("K89").Value = 3371200
$ws.Range("L89").Value = 1173800
$ws.Range("M89").Value = 828700
# Row 91
$ws.Range("D91").Value = -350100
$ws.Range("E91").Value = -368000
$ws.Range("F91").Value = -438700
$ws.Range("G91").Value = -479100
$ws.Range("H91").Value = -348900
$ws.Range("I91").Value = -397200
$ws.Range("J91").Value = -357900
$ws.Range("K91").Value = -341600
$ws.Range("L91").Value = -393200
$ws.Range("M91").Value = -305200
# Row 92
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
# Row 93
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
# Row 94
$ws.Range("D94").Value = -512700
$ws.Range("E94").Value = -427500
$ws.Range("F94").Value = -2085800
$ws.Range("G94").Value = -414000
$ws.Range("H94").Value = -534100
$ws.Range("I94").Value = -304100
$ws.Range("J94").Value = 313000
$ws.Range("K94").Value = -755900
$ws.Range("L94").Value = -1204300
$ws.Range("M94").Value = -436700
# Row 96
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
# Row 97
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
# Row 98
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
# Row 99
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
# Row 100
$ws.Range("D100").Value = 4210800
$ws.Range("E100").Value = -164900
$ws.Range("F100").Value = -2025200
$ws.Range("G100").Value = 1662800
$ws.Range("H100").Value = -673200
$ws.Range("I100").Value = -336600
$ws.Range("J100").Value = -2812800
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = -422600
$ws.Range("M100").Value = -520000
# Row 101
$ws.Range("D101").Value = 6700
$ws.Range("E101").Value = 23600
$ws.Range("F101").Value = 194100
$ws.Range("G101").Value = -115600
$ws.Range("H101").Value = -59500
$ws.Range("I101").Value = -64000
$ws.Range("J101").Value = -126800
$ws.Range("K101").Value = 5900
$ws.Range("L101").Value = -29300
$ws.Range("M101").Value = 17600
# Row 102
$ws.Range("D102").Value = 4623700
$ws.Range("E102").Value = -9000
$ws.Range("F102").Value = -3459100
$ws.Range("G102").Value = 4024600
$ws.Range("H102").Value = -234500
$ws.Range("I102").Value = -18000
$ws.Range("J102").Value = -1907400
$ws.Range("K102").Value = 2622300
$ws.Range("L102").Value = -482400
$ws.Range("M102").Value = -110300